# Weekly update: a new day's price record for Cilantro at
# "Terminal La Palmera de La Serena" is inserted as the new row 66,
# shifting the existing rows 66-97 down to 67-98 (dimension grows to R98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66 (pushes old rows 66..97 down to 67..98).
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new weekly record.
$ws.Range("A66").Value2 = 8
$ws.Range("B66").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C66").Value2 = "Coquimbo"
$ws.Range("D66").Value2 = 44523
$ws.Range("E66").Value2 = 4
$ws.Range("F66").Value2 = 100112040
$ws.Range("G66").Value2 = "Cilantro"
$ws.Range("H66").Value2 = "Sin especificar"
$ws.Range("I66").Value2 = "Primera"
$ws.Range("J66").Value2 = 3280
$ws.Range("K66").Value2 = 1500
$ws.Range("L66").Value2 = 2000
$ws.Range("M66").Value2 = 1750
$ws.Range("N66").Value2 = "$/atado 1 a 1,5 kilos"
$ws.Range("O66").Value2 = "Provincia del Elquí"
$ws.Range("P66").Value2 = 1167
$ws.Range("Q66").Value2 = 1.5
$ws.Range("R66").Value2 = "Hortaliza"
